# "calculated rank and box office $ correlations"
#
# Adds a new "Correlation" worksheet (after "Weekly Data") that pulls the
# per-week Rank / Box-Office-Gross pairs (for weeks that actually had a
# rank) out of the "Weekly Data" sheet, grouped by movie, and computes the
# CORREL() of rank vs gross overall and per-movie.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new sheet at the end of the tab strip and name it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Correlation"

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "Box Office Gross"
$ws.Range("C1").Value = "Movie"

# ---------------------------------------------------------------------
# 3. Data rows: Rank / Gross / Movie, pulled from the weeks in
#    "Weekly Data" that have a numeric Rank (i.e. excluding "No list"
#    weeks), grouped Half-Blood Prince, then Deathly Hallows Pt. 1, then
#    Deathly Hallows Pt. 2.
# ---------------------------------------------------------------------
$data = @(
    @(2, 136241423, "Half-Blood Prince"),
    @(2, 76119425, "Half-Blood Prince"),
    @(2, 37532997, "Half-Blood Prince"),
    @(3, 21083840, "Half-Blood Prince"),
    @(4, 11260386, "Half-Blood Prince"),
    @(6, 6929882, "Half-Blood Prince"),
    @(8, 4385748, "Half-Blood Prince"),
    @(9, 2643730, "Half-Blood Prince"),
    @(9, 12867, "Half-Blood Prince"),
    @(8, 3550, "Half-Blood Prince"),
    @(5, 125017372, "Deathly Hallows Pt. 1"),
    @(3, 94038757, "Deathly Hallows Pt. 1"),
    @(3, 25462987, "Deathly Hallows Pt. 1"),
    @(4, 13153605, "Deathly Hallows Pt. 1"),
    @(3, 8050897, "Deathly Hallows Pt. 1"),
    @(4, 7271017, "Deathly Hallows Pt. 1"),
    @(6, 10538580, "Deathly Hallows Pt. 1"),
    @(7, 4128555, "Deathly Hallows Pt. 1"),
    @(9, 2141702, "Deathly Hallows Pt. 1"),
    @(9, 1537296, "Deathly Hallows Pt. 1"),
    @(7, 882724, "Deathly Hallows Pt. 1"),
    @(8, 458705, "Deathly Hallows Pt. 1"),
    @(10, 532237, "Deathly Hallows Pt. 1"),
    @(10, 255162, "Deathly Hallows Pt. 1"),
    @(9, 171370, "Deathly Hallows Pt. 1"),
    @(10, 121799, "Deathly Hallows Pt. 1"),
    @(2, 133485465, "Deathly Hallows Pt. 2"),
    @(2, 125449251, "Deathly Hallows Pt. 2"),
    @(2, 52614757, "Deathly Hallows Pt. 2"),
    @(2, 27300651, "Deathly Hallows Pt. 2"),
    @(2, 15857926, "Deathly Hallows Pt. 2"),
    @(3, 9752263, "Deathly Hallows Pt. 2"),
    @(3, 5498700, "Deathly Hallows Pt. 2"),
    @(5, 3656003, "Deathly Hallows Pt. 2"),
    @(5, 3202184, "Deathly Hallows Pt. 2"),
    @(5, 1157544, "Deathly Hallows Pt. 2"),
    @(7, 734299, "Deathly Hallows Pt. 2"),
    @(9, 458141, "Deathly Hallows Pt. 2"),
    @(10, 459344, "Deathly Hallows Pt. 2"),
    @(10, 528882, "Deathly Hallows Pt. 2"),
    @(8, 100538, "Deathly Hallows Pt. 2"),
    @(4, 53011, "Deathly Hallows Pt. 2"),
    @(8, 18467, "Deathly Hallows Pt. 2")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# ---------------------------------------------------------------------
# 4. Correlation labels + formulas in column E.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Overall Correlation:"
$ws.Range("E2").Formula = "=CORREL(A2:A44,B2:B44)"

$ws.Range("E4").Value = "Half-Blood Prince Correlation:"
$ws.Range("E5").Formula = "=CORREL(A2:A11,B2:B11)"

$ws.Range("E7").Value = "Deathly Hallows 1 Correlation:"
$ws.Range("E8").Formula = "=CORREL(A12:A27,B12:B27)"

$ws.Range("E10").Value = "Deathly Hallows 2 Correlation:"
$ws.Range("E11").Formula = "=CORREL(A28:A44,B28:B44)"

# ---------------------------------------------------------------------
# 5. Column widths (best-fit-like) for B, C, E.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 17.83
$ws.Columns.Item(5).ColumnWidth = 25.33

# ---------------------------------------------------------------------
# 6. View state: the new sheet is the active tab, selection at H5.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("H5").Select()

# ---------------------------------------------------------------------
# 7. Update the view state left on the other two re-scrolled sheets.
# ---------------------------------------------------------------------
$wsBoxOffice = $wb.Worksheets.Item("Box Office")
$wsBoxOffice.Activate()
$wsBoxOffice.Range("C158").Select()

$wsWeekly = $wb.Worksheets.Item("Weekly Data")
$wsWeekly.Activate()
$wsWeekly.Range("B84").Select()

# Leave "Correlation" as the active sheet/tab.
$ws.Activate()
$ws.Range("H5").Select()
